# Predicer input_data.xlsx -- "Corrections and output methods"
#
# 1) nodes: add a "state_loss" column (K) with a non-zero value for hp1,
#    and bump hp1's initial_state (J4) from 0 to 10.
# 2) Add a new "cap_ts" time-series sheet (hp1,elc capacity profile for s1/s2/s3).
# 3) Misc selection / active-tab bookkeeping left over from the editing session.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) nodes sheet: initial_state correction (J4 0 -> 10)
# ---------------------------------------------------------------------------
$nodes = $wb.Worksheets("nodes")

$nodes.Range("J4").Value = 10

# ---------------------------------------------------------------------------
# 2) New cap_ts sheet (created/populated first so the shared-string table
#    picks up "hp1,elc,s1/s2/s3" before nodes!K1's "state_loss")
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$capTs = $wb.Worksheets.Add($null, $lastSheet)
$capTs.Name = "cap_ts"

$capTs.Range("A1").Value = "t"
$capTs.Range("B1").Value = "hp1,elc,s1"
$capTs.Range("C1").Value = "hp1,elc,s2"
$capTs.Range("D1").Value = "hp1,elc,s3"

$times = @(
    0,
    0.04166666666666669905,
    0.08333333333333330095,
    0.125,
    0.16666666666666699,
    0.20833333333333301,
    0.25,
    0.29166666666666702,
    0.33333333333333298,
    0.375,
    0.41666666666666702,
    0.45833333333333298,
    0.5,
    0.54166666666666696,
    0.58333333333333304,
    0.625,
    0.66666666666666696,
    0.70833333333333304,
    0.75,
    0.79166666666666696,
    0.83333333333333304,
    0.875,
    0.91666666666666696,
    0.95833333333333304
)

$values = @(
    5,
    4.2857142857142856,
    4.2857142857142856,
    4.4117647058823533,
    4.5454545454545459,
    4.6875,
    4.838709677419355,
    5,
    6,
    6,
    6,
    6,
    6,
    6,
    6,
    6,
    6,
    3.75,
    3.75,
    3.75,
    3.75,
    4.2857142857142856,
    4.2857142857142856,
    4.2857142857142856
)

for ($i = 0; $i -lt 24; $i++) {
    $r = $i + 2
    $capTs.Cells.Item($r, 1).Value = $times[$i]
    $capTs.Cells.Item($r, 2).Value = $values[$i]
    $capTs.Cells.Item($r, 3).Value = $values[$i]
    $capTs.Cells.Item($r, 4).Value = $values[$i]
}

$capTs.Range("A2:A25").NumberFormat = "h:mm"
$capTs.Range("M29").Select()

# ---------------------------------------------------------------------------
# 3) nodes sheet: new state_loss column (written after cap_ts so the
#    shared-string table order matches: hp1,elc,s1/s2/s3 then state_loss)
# ---------------------------------------------------------------------------
$nodes.Range("K1").Value = "state_loss"
$nodes.Range("K2").Value = 0
$nodes.Range("K3").Value = 0
$nodes.Range("K4").Value = 0.001
$nodes.Range("K5").Value = 0
$nodes.Range("K6").Value = 0
$nodes.Range("K7").Value = 0

# Copy the formatting used by the rest of the "numeric data" columns (style s="8")
$nodes.Range("J1").Copy()
$nodes.Range("K1").PasteSpecial(-4122)
$nodes.Range("J2").Copy()
$nodes.Range("K2:K3").PasteSpecial(-4122)
$nodes.Range("J5:J7").Copy()
$nodes.Range("K5:K7").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Selections / active sheet bookkeeping
# ---------------------------------------------------------------------------
$wb.Worksheets("processes").Activate()
$wb.Worksheets("processes").Range("I13").Select()

$wb.Worksheets("gen_constraint").Activate()
$wb.Worksheets("gen_constraint").Range("A2:A25").Select()

$nodes.Activate()
$nodes.Range("J5").Select()
